$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the LastName/FirstName columns (B <-> C) for the header row and all 10 data rows
for ($r = 1; $r -le 11; $r++) {
    $bVal = $ws.Cells.Item($r, 2).Value2
    $cVal = $ws.Cells.Item($r, 3).Value2
    $ws.Cells.Item($r, 2).Value2 = $cVal
    $ws.Cells.Item($r, 3).Value2 = $bVal
}

# Update the first CNE value; the rest of the column (A3:A11) recalculates via its formula
$ws.Range("A2").Value2 = 20000001

# Update the active selection to match the edited workbook
$ws.Range("G9").Select()
